$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2 (pushes existing data rows 2..21 down to 3..21)
$ws.Rows.Item(2).Insert()
# The insert operation copies the header row's bold/centered formatting onto
# the new row; strip it so the new data row matches the plain (unstyled) data
# rows elsewhere in the sheet.
$ws.Rows.Item(2).ClearFormats()

# Fill the newly inserted row 2 with its data
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "struggle"
$ws.Cells.Item(2,3).Value = 0.0610952377319335
$ws.Cells.Item(2,4).Value = 0.5822855234146118
$ws.Cells.Item(2,5).Value = 0.008422106504440301
$ws.Cells.Item(2,6).Value = -0.03398093824483917
$ws.Cells.Item(2,7).Value = -2.852332382786031
$ws.Cells.Item(2,8).Value = 0.3628882449500415

# Append 9 brand-new rows of data at the bottom (rows 23..31)
$ws.Cells.Item(23,1).Value = 2100
$ws.Cells.Item(23,2).Value = "struggle"
$ws.Cells.Item(23,3).Value = -3.810809135437012
$ws.Cells.Item(23,4).Value = 1.403007388114929
$ws.Cells.Item(23,5).Value = 0.0495486259460449
$ws.Cells.Item(23,6).Value = -1.581159264457474
$ws.Cells.Item(23,7).Value = -3.781517471585969
$ws.Cells.Item(23,8).Value = -2.15734222470499

$ws.Cells.Item(24,1).Value = 2200
$ws.Cells.Item(24,2).Value = "struggle"
$ws.Cells.Item(24,3).Value = -1.585423946380615
$ws.Cells.Item(24,4).Value = 2.060841083526612
$ws.Cells.Item(24,5).Value = -2.507726192474365
$ws.Cells.Item(24,6).Value = -0.5357818153439736
$ws.Cells.Item(24,7).Value = -0.6552340047700065
$ws.Cells.Item(24,8).Value = 0.5789350879435637

$ws.Cells.Item(25,1).Value = 2300
$ws.Cells.Item(25,2).Value = "struggle"
$ws.Cells.Item(25,3).Value = -5.486822128295898
$ws.Cells.Item(25,4).Value = 2.457437515258789
$ws.Cells.Item(25,5).Value = -1.076503276824951
$ws.Cells.Item(25,6).Value = -0.4312272305999488
$ws.Cells.Item(25,7).Value = -0.03695735122476339
$ws.Cells.Item(25,8).Value = -0.2086323031357354

$ws.Cells.Item(26,1).Value = 2400
$ws.Cells.Item(26,2).Value = "struggle"
$ws.Cells.Item(26,3).Value = 3.813155174255371
$ws.Cells.Item(26,4).Value = -5.157403945922852
$ws.Cells.Item(26,5).Value = 7.194998264312744
$ws.Cells.Item(26,6).Value = 0.1521366113910867
$ws.Cells.Item(26,7).Value = 0.3846518628451288
$ws.Cells.Item(26,8).Value = -0.3579327458021597

$ws.Cells.Item(27,1).Value = 2500
$ws.Cells.Item(27,2).Value = "struggle"
$ws.Cells.Item(27,3).Value = -3.507768154144287
$ws.Cells.Item(27,4).Value = 2.501498937606812
$ws.Cells.Item(27,5).Value = 0.7795240879058838
$ws.Cells.Item(27,6).Value = 0.3417635331956719
$ws.Cells.Item(27,7).Value = 0.5353018106246487
$ws.Cells.Item(27,8).Value = 0.2242374224018085

$ws.Cells.Item(28,1).Value = 2600
$ws.Cells.Item(28,2).Value = "struggle"
$ws.Cells.Item(28,3).Value = 0.2215757369995117
$ws.Cells.Item(28,4).Value = -0.4009582996368408
$ws.Cells.Item(28,5).Value = 2.163901329040528
$ws.Cells.Item(28,6).Value = -0.1421539567563001
$ws.Cells.Item(28,7).Value = 0.2858568746216445
$ws.Cells.Item(28,8).Value = -0.08254160519157175

$ws.Cells.Item(29,1).Value = 2700
$ws.Cells.Item(29,2).Value = "struggle"
$ws.Cells.Item(29,3).Value = 0.1625576019287109
$ws.Cells.Item(29,4).Value = 1.34720504283905
$ws.Cells.Item(29,5).Value = -0.6319388151168823
$ws.Cells.Item(29,6).Value = -0.06803667803808118
$ws.Cells.Item(29,7).Value = 0.1790894811250734
$ws.Cells.Item(29,8).Value = 0.08434615633925557

$ws.Cells.Item(30,1).Value = 2800
$ws.Cells.Item(30,2).Value = "struggle"
$ws.Cells.Item(30,3).Value = 0.044438362121582
$ws.Cells.Item(30,4).Value = -0.1398162841796875
$ws.Cells.Item(30,5).Value = -0.8414495587348938
$ws.Cells.Item(30,6).Value = -0.08848196070413178
$ws.Cells.Item(30,7).Value = -0.1032362286837731
$ws.Cells.Item(30,8).Value = 0.232543302129726

$ws.Cells.Item(31,1).Value = 2900
$ws.Cells.Item(31,2).Value = "struggle"
$ws.Cells.Item(31,3).Value = -0.1983919143676757
$ws.Cells.Item(31,4).Value = -0.413076639175415
$ws.Cells.Item(31,5).Value = 0.2017757892608642
$ws.Cells.Item(31,6).Value = 0.06249837318853448
$ws.Cells.Item(31,7).Value = 0.04699299066346516
$ws.Cells.Item(31,8).Value = 0.2259266389419836

# The "timestamp" column (A) is a simple row-based sequence (0, 100, 200, ...)
# and must NOT shift along with the rest of the row data when the new row is
# inserted, so re-normalize it for every data row now that the sheet has 30
# data rows (rows 2..31).
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r,1).Value = ($r - 2) * 100
}
